# Generate Report for Handback
# Update file GUIDs / hashes / timestamps produced by a new handback run.
#
# Sheet "Overview" (index 1): column A file-name links
# Sheet "zh-cn" (index 2): Source/Target file-name + xlf links, handoff/handback datetimes
# Sheet "de-de" (index 3): Source/Target file-name + xlf links, handoff/handback datetimes

$wb = $excel.ActiveWorkbook

$newMd1  = "69d0169b-6270-4df0-84d7-709428d9db77.md"
$newMd2  = "ffff4574e5b5-b010-44d7-b42e-505b004d5a32.md"
$newXlfZh = "69d0169b-6270-4df0-84d7-709428d9db77.ffa7541d8c2d17af2c9da6597237242db4946ebb.zh-cn.xlf"
$newXlfDe = "69d0169b-6270-4df0-84d7-709428d9db77.ffa7541d8c2d17af2c9da6597237242db4946ebb.de-de.xlf"

$dtZh1 = "2016-03-23 15:13:56"
$dtZh2 = "2016-03-23 15:14:20"
$dtDe1 = "2016-03-23 15:14:01"
$dtDe2 = "2016-03-23 15:14:27"

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item(1)

# Preserve the original (external) hyperlink targets - only the displayed
# text / underlying cell text actually changes with this edit.
$overviewLinkA2 = "https://github.com/OpenLocalizationTest/oltest/blob/f3df70e447d057b39401f5fa56f6afe1a0449433/e2e/45ebf1a3-926a-4c76-85b5-dfa0cb46f77a.md"
$overviewLinkA3 = "https://github.com/OpenLocalizationTest/oltest/blob/f3df70e447d057b39401f5fa56f6afe1a0449433/e2e/988ddba3-a78c-45d3-8b3d-e4fc07f14f5d.md"

# Clear existing hyperlinks on the sheet so they can be re-created with the
# updated display text (the API only supports adding hyperlinks, not
# editing the display text of an existing one in place).
$wsOverview.Range("A2").Hyperlinks.Delete()

$wsOverview.Range("A2").Value = $newMd1
$wsOverview.Range("A3").Value = $newMd2

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $overviewLinkA2, "", "", $newMd1)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $overviewLinkA3, "", "", $newMd2)

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item(2)

$zhLinks = @{
    "A2" = "https://github.com/OpenLocalizationTest/oltest/blob/f3df70e447d057b39401f5fa56f6afe1a0449433/e2e/45ebf1a3-926a-4c76-85b5-dfa0cb46f77a.md"
    "D2" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cc080092d3126f39559f47c8a59942db5d7bea36/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/45ebf1a3-926a-4c76-85b5-dfa0cb46f77a.d4b0174df86b94ece4aaaa7e79f7bfd89845e8ab.zh-cn.xlf"
    "F2" = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/a459cd2c619c5c370537b39c7652a31c39f4fc61/e2e/45ebf1a3-926a-4c76-85b5-dfa0cb46f77a.md"
    "G2" = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6584873f15028d47ca01b7389d782c9c042146d3/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/45ebf1a3-926a-4c76-85b5-dfa0cb46f77a.d4b0174df86b94ece4aaaa7e79f7bfd89845e8ab.zh-cn.xlf"
    "A3" = "https://github.com/OpenLocalizationTest/oltest/blob/f3df70e447d057b39401f5fa56f6afe1a0449433/e2e/988ddba3-a78c-45d3-8b3d-e4fc07f14f5d.md"
    "D3" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cc080092d3126f39559f47c8a59942db5d7bea36/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/988ddba3-a78c-45d3-8b3d-e4fc07f14f5d.4c99c016559bbe1e9781807c55abe70de50afe42.zh-cn.xlf"
    "F3" = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/a459cd2c619c5c370537b39c7652a31c39f4fc61/e2e/988ddba3-a78c-45d3-8b3d-e4fc07f14f5d.md"
    "G3" = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6584873f15028d47ca01b7389d782c9c042146d3/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/988ddba3-a78c-45d3-8b3d-e4fc07f14f5d.4c99c016559bbe1e9781807c55abe70de50afe42.zh-cn.xlf"
}

$zhValues = @{
    "A2" = $newMd1
    "F2" = $newMd1
    "A3" = $newMd2
    "F3" = $newMd2
    "D2" = $newXlfZh
    "G2" = $newXlfZh
    "D3" = $newXlfZh
    "G3" = $newXlfZh
}

$wsZh.Range("A2").Hyperlinks.Delete()

foreach ($addr in $zhValues.Keys) {
    $wsZh.Range($addr).Value = $zhValues[$addr]
}
$wsZh.Range("E2").Value = $dtZh1
$wsZh.Range("E3").Value = $dtZh1
$wsZh.Range("H2").Value = $dtZh2
$wsZh.Range("H3").Value = $dtZh2

foreach ($addr in @("A2","D2","F2","G2","A3","D3","F3","G3")) {
    $wsZh.Hyperlinks.Add($wsZh.Range($addr), $zhLinks[$addr], "", "", $zhValues[$addr])
}

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item(3)

$deLinks = @{
    "A2" = "https://github.com/OpenLocalizationTest/oltest/blob/f3df70e447d057b39401f5fa56f6afe1a0449433/e2e/45ebf1a3-926a-4c76-85b5-dfa0cb46f77a.md"
    "D2" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2b56159cc6508967bfc701e21678d612e8f7c69e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/45ebf1a3-926a-4c76-85b5-dfa0cb46f77a.d4b0174df86b94ece4aaaa7e79f7bfd89845e8ab.de-de.xlf"
    "F2" = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/d7d7f02d1b1175ceb9cc1520ded175f83c9f2a64/e2e/45ebf1a3-926a-4c76-85b5-dfa0cb46f77a.md"
    "G2" = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/46a39ad94d4a628dff86a72d3f68e006d2ba4786/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/45ebf1a3-926a-4c76-85b5-dfa0cb46f77a.d4b0174df86b94ece4aaaa7e79f7bfd89845e8ab.de-de.xlf"
    "A3" = "https://github.com/OpenLocalizationTest/oltest/blob/f3df70e447d057b39401f5fa56f6afe1a0449433/e2e/988ddba3-a78c-45d3-8b3d-e4fc07f14f5d.md"
    "D3" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2b56159cc6508967bfc701e21678d612e8f7c69e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/988ddba3-a78c-45d3-8b3d-e4fc07f14f5d.4c99c016559bbe1e9781807c55abe70de50afe42.de-de.xlf"
    "F3" = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/d7d7f02d1b1175ceb9cc1520ded175f83c9f2a64/e2e/988ddba3-a78c-45d3-8b3d-e4fc07f14f5d.md"
    "G3" = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/46a39ad94d4a628dff86a72d3f68e006d2ba4786/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/988ddba3-a78c-45d3-8b3d-e4fc07f14f5d.4c99c016559bbe1e9781807c55abe70de50afe42.de-de.xlf"
}

$deValues = @{
    "A2" = $newMd1
    "F2" = $newMd1
    "A3" = $newMd2
    "F3" = $newMd2
    "D2" = $newXlfDe
    "G2" = $newXlfDe
    "D3" = $newXlfDe
    "G3" = $newXlfDe
}

$wsDe.Range("A2").Hyperlinks.Delete()

foreach ($addr in $deValues.Keys) {
    $wsDe.Range($addr).Value = $deValues[$addr]
}
$wsDe.Range("E2").Value = $dtDe1
$wsDe.Range("E3").Value = $dtDe1
$wsDe.Range("H2").Value = $dtDe2
$wsDe.Range("H3").Value = $dtDe2

foreach ($addr in @("A2","D2","F2","G2","A3","D3","F3","G3")) {
    $wsDe.Hyperlinks.Add($wsDe.Range($addr), $deLinks[$addr], "", "", $deValues[$addr])
}
